$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.530583262443542
$ws.Range("B1").Value = 2.306761264801025
$ws.Range("C1").Value = 4.234432697296143
$ws.Range("D1").Value = 1.853784441947937
$ws.Range("E1").Value = 0.8239727020263672
